$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TasksDone")

# --- Tasks 21-30 block: mark tasks T21-T30 as done (2 points each) ---
$ws.Range("C35").Value = 2
$ws.Range("C36").Value = 2
$ws.Range("C37").Value = 2
$ws.Range("C38").Value = 2
$ws.Range("C39").Value = 2
$ws.Range("C40").Value = 2
$ws.Range("C41").Value = 2
$ws.Range("C42").Value = 2
$ws.Range("C43").Value = 2
$ws.Range("C44").Value = 2

# --- Tasks 31-42 block: mark tasks T31-T40 as done (2 points each), T41/T42 at 3 ---
$ws.Range("C48").Value = 2
$ws.Range("C49").Value = 2
$ws.Range("C50").Value = 2
$ws.Range("C51").Value = 2
$ws.Range("C52").Value = 2
$ws.Range("C53").Value = 2
$ws.Range("C54").Value = 2
$ws.Range("C55").Value = 2
$ws.Range("C56").Value = 2
$ws.Range("C57").Value = 3

# Insert a new row right before the current T42 row (row 59), which pushes
# the existing T42 row (and its bottom border formatting) down to row 60.
$ws.Rows("59").Insert()

# Row 59 is now a fresh blank row above the (shifted) T42 row: re-fill it
# with the T42 data so T42 keeps its original position/content.
$t42Label = $ws.Range("B60").Value2
$ws.Range("A59").Value = 42
$ws.Range("B59").Value = $t42Label
$ws.Range("C59").Value = 3

# Row 60 (which still carries the thick-bottom border that used to mark the
# end of the 31-42 block) becomes the new bonus task row.
$ws.Range("A60").Value = 43
$ws.Range("B60").Value = "T43 Bonus Task"
$ws.Range("C60").Value = 3

# Update the "sum of tasks 31-42" formula to include the new row.
$ws.Range("C61").Formula = "=SUM(C48:C60)"

# Update the grand-total formula to point at the (shifted) sum row.
$ws.Range("C63").Formula = "=`$C`$19+`$C`$32+`$C`$45+`$C`$61"

# Reflect the scrolled/selected view state from the saved workbook.
[void]$ws.Range("F61").Select()
$excel.ActiveWindow.ScrollRow = 53
$excel.ActiveWindow.ScrollColumn = 1
